$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (data index 0-23),
# columns B,C,D,E,F,K,M,N,O (bus-related loading percentages).
$columns = @("B","C","D","E","F","K","M","N","O")

$data = @(
    @(9.07207546586166, 5.459457105089298, 5.030686024327975, 12.60419834982519, 24.92074885364291, 8.318395326212755, 13.55570814820479, 18.73536573970178, 22.30085422391774),
    @(8.793952934427979, 5.337939148169157, 4.980291331590418, 12.38848904863236, 24.92046979955405, 8.131818729482857, 13.38679700773052, 18.79653686691411, 22.34993385733469),
    @(8.620428974443081, 5.261101922172143, 4.948568849471397, 12.25847303576059, 24.92681840777126, 8.016568215339843, 13.28538495281228, 18.8358296404732, 22.38488491961896),
    @(8.549139841634814, 5.229256138853356, 4.935451079583654, 12.20617840493596, 24.93104492068234, 7.969497817851952, 13.24468607267562, 18.85227894661533, 22.40033627026948),
    @(8.537271056031736, 5.223936691972641, 4.933261535356899, 12.19753869556962, 24.93184572232512, 7.961677542968792, 13.23796735225356, 18.85503678404612, 22.4029748684064),
    @(8.619469719911105, 5.260674566191156, 4.948392703019005, 12.25776488587958, 24.92686877105337, 8.015933739199028, 13.28483346955454, 18.83604970973238, 22.38508841223688),
    @(8.976817951664962, 5.418034233351058, 5.013475085774155, 12.52936370980556, 24.91929967815921, 8.254247970695754, 13.49701995160021, 18.75609854685852, 22.31677588811801),
    @(9.651025540820896, 5.707926457291405, 5.134677680861465, 13.07784817311722, 24.95614923766896, 8.713189693693817, 13.92916458955616, 18.61300869739646, 22.22113129400537),
    @(10.12450650574261, 5.908240153230784, 5.219504811947827, 13.48596245894423, 25.01461590410904, 9.041585856328888, 14.25345350989057, 18.51614322706785, 22.17434832336479),
    @(10.33413131093121, 5.996378294839179, 5.25711664657133, 13.67180371275984, 25.04798163833378, 9.188365516575523, 14.40181944346933, 18.47385246489195, 22.15818765175469),
    @(10.41260754905408, 6.029306675219007, 5.271213807230007, 13.74212883153385, 25.06158364639455, 9.243518024647193, 14.45807078135601, 18.45809177892035, 22.15280573585068),
    @(10.39574761277598, 6.02223514278726, 5.268184289387827, 13.7269864135889, 25.05861130896162, 9.231659885792627, 14.44595384402791, 18.46147485047707, 22.15393199634923),
    @(10.34060612129976, 5.999096402968362, 5.25827937627044, 13.67759080802451, 25.04908134088673, 9.192911858170346, 14.40644620549688, 18.47255074172407, 22.15773008537397),
    @(10.30671052323028, 5.984864445866256, 5.252193214453348, 13.64732599594875, 25.04336970398213, 9.169120052767104, 14.38225397980386, 18.47936807678274, 22.16015264049465),
    @(10.11068442850401, 5.902418442252207, 5.217026692101895, 13.47381452064807, 25.01257107683991, 9.03193620148881, 14.24376988968277, 18.51894262182155, 22.17550763552607),
    @(9.988896360876478, 5.851062311722032, 5.195199316796296, 13.36736407085958, 24.99540673929869, 8.947069981301196, 14.15899276365411, 18.54367380973269, 22.18624014414038),
    @(9.918308523074671, 5.821243604732286, 5.182553181078505, 13.30615866513676, 24.98617183988144, 8.898014581437332, 14.11031389705276, 18.55806556482413, 22.19289516244541),
    @(9.894318619669633, 5.81109999183481, 5.178255844751438, 13.2854417190357, 24.98315473653653, 8.881365325745373, 14.09384790722108, 18.56296708908605, 22.19523116825381),
    @(10.00191721818537, 5.856558377772313, 5.197532397380421, 13.37869417740265, 24.99716796734013, 8.956129650392867, 14.16800926310409, 18.54102385405698, 22.18504775798563),
    @(10.35682761232696, 6.005905098423753, 5.261192682288145, 13.69210141368768, 25.0518543314005, 9.204305172023561, 14.41804911964118, 18.46929060223738, 22.15659446208823),
    @(10.58348270643638, 6.100895299634944, 5.301947245523516, 13.89661285547321, 25.09322848717158, 9.363974696490091, 14.58183928465978, 18.42388821082192, 22.14229934475539),
    @(10.46302035129787, 6.050442311631231, 5.280275324058471, 13.78751470104583, 25.07063322034969, 9.279004340172301, 14.49440440430095, 18.44798533578674, 22.14953502809094),
    @(9.996032259007631, 5.854074518757951, 5.196477914273249, 13.37357184794317, 24.99636974358046, 8.95203459421834, 14.16393271076902, 18.54222135857185, 22.18558532586309),
    @(9.472096331963353, 5.631641703196576, 5.102608798974005, 12.92824092270028, 24.94065427313301, 8.590324488098165, 13.81084798179818, 18.6502611353412, 22.24288933889782)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range($columns[$j] + $row).Value = $rowValues[$j]
    }
}
